$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new data row (item 3: "FUSI-ZON CREAM 15 GM") above the
#     current row 6 ("MOTINORM ..."), pushing everything below down by one.
$ws.Range("A6:N6").Insert()

# Copy the formatting (styles) of the row that now sits at row 7 (the old
# row 6 / MOTINORM row) onto the freshly inserted blank row 6, so the new
# row reuses the existing table styles instead of generating new ones.
$ws.Range("A7:N7").Copy()
$ws.Range("A6:N6").PasteSpecial(-4122)
$ws.Rows.Item(6).RowHeight = 24.75

# Recreate the merged cells for the new row (B:G, H:K, L:M), matching the
# pattern used by every other data row in the table.
$ws.Range("B6:G6").Merge()
$ws.Range("H6:K6").Merge()
$ws.Range("L6:M6").Merge()

# Populate the new row's data.
$ws.Range("A6").Value = 3
$ws.Range("B6").Value = "FUSI-ZON CREAM 15 GM"
$ws.Range("H6").Value = "1:0"
$ws.Range("L6").Value = 48
$ws.Range("N6").Value = "1:0"

# Renumber the "م" (item index) column for the rows that were pushed down.
$ws.Range("A7").Value = 4
$ws.Range("A8").Value = 5
$ws.Range("A9").Value = 6
$ws.Range("A10").Value = 7

# --- The totals row (previously row 10, now shifted to row 11) needs its
#     grand total bumped by the new row's price.
$ws.Range("K11").Value = 264.36

# Match the refreshed layout's row heights.
$ws.Rows.Item(7).RowHeight = 25.5
$ws.Rows.Item(8).RowHeight = 25.5
$ws.Rows.Item(9).RowHeight = 24.75
$ws.Rows.Item(10).RowHeight = 25.5
$ws.Rows.Item(11).RowHeight = 25.5
